# Update workbook: rename activity sheets and refresh header/footer font style.
$wb = $excel.ActiveWorkbook

# --- Rename the two worksheets (Activité 8/9 -> Activité 12/13) ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws1.Name = "Activité 12"
$ws2.Name = "Activité 13"

# --- Refresh header/footer font style name on both sheets:
#     "Times New Roman,Regular" -> "Times New Roman,Normal" ---
$ps1 = $ws1.PageSetup
$ps1.CenterHeader = '&"Times New Roman,Normal"&12&A'
$ps1.CenterFooter = '&"Times New Roman,Normal"&12Page &P'
# Re-assert the fit-to-page settings so they are preserved on save.
$ps1.FitToPagesWide = 1
$ps1.FitToPagesTall = 1

$ps2 = $ws2.PageSetup
$ps2.CenterHeader = '&"Times New Roman,Normal"&12&A'
$ps2.CenterFooter = '&"Times New Roman,Normal"&12Page &P'
$ps2.FitToPagesWide = 1
$ps2.FitToPagesTall = 1
